# Commit: "improved territory_options performance (hardcoded via dict comphr)"
#
# Rename the (only) worksheet to reflect that it now holds the deprecated
# / original data tables, and move the live cell selection to C22 (the
# cell the author was working in when the sheet was renamed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet tab rename: "final data tables" -> "original data tables - deprecat"
$ws.Name = "original data tables - deprecat"

# Make sure this sheet is the active one, then move the selection.
$ws.Activate()
$ws.Range("C22").Select()
